$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: 2018-02-13, 2 hours, "R work"
$ws.Cells.Item(10, 1).Value = 43144
$ws.Cells.Item(10, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(10, 2).Value = 2
$ws.Cells.Item(10, 3).Value = "R work"

# Row 11: 2018-02-18, 3 hours, "power simulation; power function"
$ws.Cells.Item(11, 1).Value = 43149
$ws.Cells.Item(11, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(11, 2).Value = 3
$ws.Cells.Item(11, 3).Value = "power simulation; power function"

# Update selection to mirror the saved workbook view state
[void]$ws.Range("A12").Select()
